# TestRenderEmptyEnumerable.xlsx edit
#
# Commit: "Add DateTime to ExtendedPrimitiveTypes, Enum-types render in
# dynamic panel by default"
#
# The dynamic-panel header row gains a new leading "Sex" column (an
# Enum-typed property that now renders by default), pushing the existing
# "Name" / "Sum" / "Date" headers one column to the right, and the
# template's format-placeholder row shifts along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Name" header) - this shifts
# B2:D2 ("Name","Sum","Date") right to C2:E2, and the template row 3
# placeholder cells (C3, D3) right to D3, E3, preserving their
# values/styles/number formats.
$ws.Columns("B:B").Insert()

# New header cell for the inserted column: the Enum ("Sex") property.
$ws.Range("B2").Value = "Sex"

# Row 3 holds only format placeholders (no actual data - this is the
# "empty enumerable" render case). Re-assert the expected number formats
# on the shifted placeholder cells: the numeric "Sum" placeholder keeps
# its money format, and the blank "Date" placeholder keeps its date
# format.
$ws.Range("D3").NumberFormat = "#,0.00"
$ws.Range("E3").NumberFormat = "yyyy-MM-dd"

# The "Date" placeholder cell carries no value (blank template cell);
# make sure the shift didn't leave stray content behind.
$ws.Range("E3").Value = ""
